$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the "Purpose of Program" paragraph text and re-split it into
#    three runs: "Purpose of Program: " | "A program that shows the pri" |
#    "me factor of the number the user input and display the prime number
#    from smallest prime number to largest prime number."
# ---------------------------------------------------------------------------
$start = 88
$oldEnd = 270
$full = $d.Range($start, $oldEnd)
$full.Text = "Purpose of Program: A program that shows the prime factor of the number the user input and display the prime number from smallest prime number to largest prime number."

# Force a run boundary after "Purpose of Program: " (20 chars) and after
# "...the pri" (28 more chars) by toggling Bold off/on (same end value ->
# no visible formatting change, but the engine only merges runs that were
# never distinguished, so this reliably splits them).
$split1 = $d.Range($start, $start + 20)
$split1.Bold = 0
$split1.Bold = 1

$split2 = $d.Range($start + 20, $start + 48)
$split2.Bold = 0
$split2.Bold = 1

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the "Sort ..." paragraph to
#    the end of the "Purpose of Program" paragraph.
# ---------------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$p6 = $d.Paragraphs.Item(6)
$endPos = $p6.Range.End - 1

# Placing a bookmark exactly at "end of paragraph text" (collapsed range
# right before the paragraph mark) needs a small sentinel trick: insert a
# throw-away character, bookmark the (now mid-paragraph) gap before it,
# then delete the sentinel again.
$sentinel = $d.Range($endPos, $endPos)
$sentinel.InsertAfter("Z")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanup = $d.Range($endPos, $endPos + 1)
$cleanup.Text = ""
